# feat: add antibiotic fields
#
# The "strains" sheet grows from a 2-column (Tag, Plasmids) lookup into a
# 4-column one (Tag, Strain, Plasmids, Antibiotics). The original "Plasmids"
# column shifts from B to C, a new "Strain" column is inserted at B, a new
# "Antibiotics" column is appended at D, and two more tag rows are added so
# each new column has a sample value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header row gains "Strain" (B) and "Antibiotics" (D); "Plasmids"
# shifts right to C.
$ws.Range("B1").Value = "Strain"
$ws.Range("C1").Value = "Plasmids"
$ws.Range("D1").Value = "Antibiotics"

# Row 2: existing tag "s2" now also gets a Strain value; the old "p2,p3"
# Plasmids value moves down to row 3 alongside its own tag.
$ws.Range("B2").Value = "s0"

# Row 3: new tag "s3" carries the Plasmids value that used to live on row 2.
$ws.Range("A3").Value = "s3"
$ws.Range("C3").Value = "p2,p3"

# Row 4: new tag "s4" carries the new Antibiotics sample value.
$ws.Range("A4").Value = "s4"
$ws.Range("D4").Value = "Amp,Kan"

# Leave the selection where the author's editing session ended up.
$ws.Range("C52").Select() | Out-Null
